# df_importance.xlsx - refresh the feature-importance table (split/gain)
# and append the new engineered "*_sent_to_contract" / "total_ether_sent_contracts"
# features that came out of the two new "Artigo - ..." notebooks.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A = feature name, B = split (int), C = gain (float).
# Rows 2-37 already existed (values refreshed below); rows 38-41 are brand new
# and need to be created first so they inherit the bold/centered/bordered
# label style used by the rest of column A (copied from the last existing row).
$newRows = 38..41
foreach ($newRow in $newRows) {
    $ws.Cells.Item(37, 1).Copy() | Out-Null
    $ws.Cells.Item($newRow, 1).PasteSpecial(-4122) | Out-Null
}

$importance = @(
    @{ Row = 2; Name = 'Avg_min_between_sent_tnx'; Split = 140; Gain = 202.248086489737 },
    @{ Row = 3; Name = 'Avg_min_between_received_tnx'; Split = 244; Gain = 1017.549922570586 },
    @{ Row = 4; Name = 'Time_Diff_between_first_and_last_(Mins)'; Split = 477; Gain = 30669.84143901244 },
    @{ Row = 5; Name = 'Sent_tnx'; Split = 136; Gain = 858.8448903486133 },
    @{ Row = 6; Name = 'Received_Tnx'; Split = 124; Gain = 960.3180639818311 },
    @{ Row = 7; Name = 'Number_of_Created_Contracts'; Split = 6; Gain = 63.79238022863865 },
    @{ Row = 8; Name = 'Unique_Received_From_Addresses'; Split = 527; Gain = 4509.413096152246 },
    @{ Row = 9; Name = 'Unique_Sent_To_Addresses'; Split = 72; Gain = 80.18373914808035 },
    @{ Row = 10; Name = 'min_value_received'; Split = 235; Gain = 437.2278335541487 },
    @{ Row = 11; Name = 'max_value_received'; Split = 156; Gain = 214.1127726882696 },
    @{ Row = 12; Name = 'avg_val_received'; Split = 363; Gain = 1715.074003368616 },
    @{ Row = 13; Name = 'min_val_sent'; Split = 191; Gain = 1713.144184119999 },
    @{ Row = 14; Name = 'max_val_sent'; Split = 53; Gain = 56.49988653510809 },
    @{ Row = 15; Name = 'avg_val_sent'; Split = 115; Gain = 184.4701158553362 },
    @{ Row = 16; Name = 'min_value_sent_to_contract'; Split = 0; Gain = 0 },
    @{ Row = 17; Name = 'max_val_sent_to_contract'; Split = 0; Gain = 0 },
    @{ Row = 18; Name = 'avg_value_sent_to_contract'; Split = 0; Gain = 0 },
    @{ Row = 19; Name = 'total_transactions_(including_tnx_to_create_contract'; Split = 363; Gain = 2534.980850402266 },
    @{ Row = 20; Name = 'total_Ether_sent'; Split = 105; Gain = 162.6454304680228 },
    @{ Row = 21; Name = 'total_ether_received'; Split = 192; Gain = 2972.544908896089 },
    @{ Row = 22; Name = 'total_ether_sent_contracts'; Split = 0; Gain = 0 },
    @{ Row = 23; Name = 'total_ether_balance'; Split = 283; Gain = 1021.243696816266 },
    @{ Row = 24; Name = 'Total_ERC20_tnxs'; Split = 413; Gain = 29276.64033755288 },
    @{ Row = 25; Name = 'ERC20_total_Ether_received'; Split = 333; Gain = 2757.246877282858 },
    @{ Row = 26; Name = 'ERC20_total_ether_sent'; Split = 125; Gain = 862.4597082287073 },
    @{ Row = 27; Name = 'ERC20_total_Ether_sent_contract'; Split = 0; Gain = 0 },
    @{ Row = 28; Name = 'ERC20_uniq_sent_addr'; Split = 59; Gain = 62.21378822624683 },
    @{ Row = 29; Name = 'ERC20_uniq_rec_addr'; Split = 102; Gain = 412.1487497240305 },
    @{ Row = 30; Name = 'ERC20_uniq_sent_addr.1'; Split = 0; Gain = 0 },
    @{ Row = 31; Name = 'ERC20_uniq_rec_contract_addr'; Split = 121; Gain = 441.3325098231435 },
    @{ Row = 32; Name = 'ERC20_min_val_rec'; Split = 203; Gain = 2112.481059491634 },
    @{ Row = 33; Name = 'ERC20_max_val_rec'; Split = 389; Gain = 2081.079403884709 },
    @{ Row = 34; Name = 'ERC20_avg_val_rec'; Split = 135; Gain = 286.1355265527964 },
    @{ Row = 35; Name = 'ERC20_min_val_sent'; Split = 91; Gain = 288.9137238487601 },
    @{ Row = 36; Name = 'ERC20_max_val_sent'; Split = 48; Gain = 63.70978651195765 },
    @{ Row = 37; Name = 'ERC20_avg_val_sent'; Split = 56; Gain = 113.7439238354564 },
    @{ Row = 38; Name = 'ERC20_uniq_sent_token_name'; Split = 75; Gain = 217.1866341531277 },
    @{ Row = 39; Name = 'ERC20_uniq_rec_token_name'; Split = 34; Gain = 50.08650804311037 },
    @{ Row = 40; Name = 'ERC20_most_sent_token_type'; Split = 5; Gain = 3.454907029867172 },
    @{ Row = 41; Name = 'ERC20_most_rec_token_type'; Split = 79; Gain = 1960.085619091988 }
)

foreach ($item in $importance) {
    $ws.Range("A$($item.Row)").Value = $item.Name
    $ws.Range("B$($item.Row)").Value = $item.Split
    $ws.Range("C$($item.Row)").Value = $item.Gain
}
